$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "52.303.33"
$ws.Range("E2").Value = "  -0.08%  "

# Row 3
$ws.Range("D3").Value = "2.852.19"
$ws.Range("E3").Value = "  +1.97%  "

# Row 4
$ws.Range("E4").Value = "  -0.01%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "362.25"
$ws.Range("E5").Value = "  +6.21%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "113.80"
$ws.Range("E6").Value = "  -3.08%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.575"
$ws.Range("E7").Value = "  +3.83%  "

# Row 8
$ws.Range("E8").Value = "  +0.04%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.603"
$ws.Range("E9").Value = "  +3.59%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "41.65"
$ws.Range("E10").Value = "  -1.47%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0864"
$ws.Range("E11").Value = "  -0.87%  "

# Row 12
$ws.Range("E12").Value = "  +1.22%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "20.01"
$ws.Range("E13").Value = "  -0.51%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.79"
$ws.Range("E14").Value = "  +1.88%  "

# Row 15
$ws.Range("D15").Value = "3.292.72"
$ws.Range("E15").Value = "  +1.90%  "

# Row 16
$ws.Range("D16").Value = "2.830.83"
$ws.Range("E16").Value = "  +1.39%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.905"
$ws.Range("E17").Value = "  +1.72%  "

# Row 18
$ws.Range("D18").Value = "52.080.17"
$ws.Range("E18").Value = "  -0.11%  "

# Row 19
$ws.Range("E19").Value = "  +8.79%  "

# Row 20
$ws.Range("E20").Value = "  -2.37%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.64"
$ws.Range("E21").Value = "  +1.61%  "

# Row 22
$ws.Range("E22").Value = "  +0.42%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "70.41"
$ws.Range("E23").Value = "  +0.00%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "268.11"
$ws.Range("E24").Value = "  -4.00%  "

# Row 25
$ws.Range("E25").Value = "  +0.09%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "27.29"
$ws.Range("E26").Value = "  +1.08%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.00"
$ws.Range("E27").Value = "  +0.06%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.45"
$ws.Range("E28").Value = "  +2.03%  "

# Row 29
$ws.Range("E29").Value = "  +1.44%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "53.60"
$ws.Range("E30").Value = "  +6.38%  "

# Row 31
$ws.Range("E31").Value = "  -1.26%  "

# Row 32
$ws.Range("B32").Value = "InjectiveProtocol"
$ws.Range("C32").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "34.28"
$ws.Range("E32").Value = "  -2.29%  "

# Row 33
$ws.Range("B33").Value = "VeChain"
$ws.Range("C33").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0457"
$ws.Range("E33").Value = "  +22.41%  "

# Row 34
$ws.Range("E34").Value = "  +2.99%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.40"
$ws.Range("E35").Value = "  +8.17%  "

# Row 36
$ws.Range("E36").Value = "  +2.05%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.998"
$ws.Range("E37").Value = "  -0.09%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.29"
$ws.Range("E38").Value = "  +0.52%  "

# Row 39
$ws.Range("E39").Value = "  -2.57%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "18.36"
$ws.Range("E40").Value = "  -3.72%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "23.82"
$ws.Range("E41").Value = "  +1.24%  "

# Row 42
$ws.Range("E42").Value = "  +1.23%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "128.96"
$ws.Range("E43").Value = "  +1.67%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.57"
$ws.Range("E44").Value = "  -6.73%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.28"
$ws.Range("E45").Value = "  -3.02%  "

# Row 46
$ws.Range("B46").Value = "NEARProtocol"
$ws.Range("C46").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.41"
$ws.Range("E46").Value = "  +1.65%  "

# Row 47
$ws.Range("B47").Value = "Maker"
$ws.Range("C47").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D47").Value = "2.118.66"
$ws.Range("E47").Value = "  +0.40%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.01"
$ws.Range("E49").Value = "  +9.92%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "5.87"
$ws.Range("E50").Value = "  +5.37%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "9.09"
$ws.Range("E51").Value = "  +1.29%  "
